$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D2:E51 so numeric-looking strings (e.g. "7.70") are not
# auto-coerced into numbers by Excel, matching the source inlineStr cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '42.695.47'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '2.545.37'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '318.90'
$ws.Range('E5').Value = '  +4.56%  '
$ws.Range('D6').Value = '95.09'
$ws.Range('E6').Value = '  -2.76%  '
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '0.535'
$ws.Range('E9').Value = '  -2.01%  '
$ws.Range('D10').Value = '36.31'
$ws.Range('E10').Value = '  -1.31%  '
$ws.Range('D11').Value = '0.0814'
$ws.Range('E11').Value = '  -1.24%  '
$ws.Range('D12').Value = '7.70'
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').Value = '2.935.57'
$ws.Range('E14').Value = '  +0.22%  '
$ws.Range('D15').Value = '15.98'
$ws.Range('E15').Value = '  +6.08%  '
$ws.Range('D16').Value = '2.510.49'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D17').Value = '0.869'
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('D18').Value = '42.739.74'
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('D19').Value = '13.07'
$ws.Range('E19').Value = '  -1.54%  '
$ws.Range('D20').Value = '6.64'
$ws.Range('E20').Value = '  +1.20%  '
$ws.Range('E21').Value = '  -1.93%  '
$ws.Range('D22').Value = '71.11'
$ws.Range('D23').Value = '252.79'
$ws.Range('E23').Value = '  -0.54%  '
$ws.Range('D24').Value = '2.98'
$ws.Range('E24').Value = '  +1.82%  '
$ws.Range('E25').Value = '  -3.01%  '
$ws.Range('D26').Value = '27.29'
$ws.Range('E26').Value = '  -2.32%  '
$ws.Range('E27').Value = '  +0.18%  '
$ws.Range('E28').Value = '  +4.19%  '
$ws.Range('D29').Value = '39.83'
$ws.Range('E29').Value = '  +4.68%  '
$ws.Range('D30').Value = '10.27'
$ws.Range('E30').Value = '  +0.97%  '
$ws.Range('D31').Value = '5.96'
$ws.Range('E31').Value = '  -3.71%  '
$ws.Range('D32').Value = '156.05'
$ws.Range('E32').Value = '  -0.75%  '
$ws.Range('D33').Value = '2.15'
$ws.Range('E33').Value = '  +0.83%  '
$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D34').Value = '19.29'
$ws.Range('E34').Value = '  -1.60%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = '3.36'
$ws.Range('E35').Value = '  +1.41%  '
$ws.Range('D36').Value = '0.0792'
$ws.Range('E36').Value = '  -0.44%  '
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('E38').Value = '  -2.64%  '
$ws.Range('D39').Value = '2.46'
$ws.Range('E39').Value = '  +14.18%  '
$ws.Range('E40').Value = '  -0.26%  '
$ws.Range('D41').Value = '23.66'
$ws.Range('E41').Value = '  -5.51%  '
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('D43').Value = '3.36'
$ws.Range('E43').Value = '  -1.64%  '
$ws.Range('E44').Value = '  +0.39%  '
$ws.Range('D45').Value = '0.0302'
$ws.Range('E45').Value = '  -0.98%  '
$ws.Range('D46').Value = '2.028.68'
$ws.Range('E46').Value = '  -2.98%  '
$ws.Range('D47').Value = '84.48'
$ws.Range('E47').Value = '  -2.44%  '
$ws.Range('E48').Value = '  -0.12%  '
$ws.Range('D49').Value = '2.789.43'
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('D50').Value = '73.82'
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('E51').Value = '  -0.55%  '

# Drop the temporary text-format style so cells keep the original (unstyled)
# look, only the stored values differ.
$ws.Range("D2:E51").ClearFormats()
